# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#  - Status column (B) moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - New columns are populated: E = Latest Target File, F = Latest Handback File
#    (both are hyperlinks, mirroring the Source File (A) / Latest Handoff File (C) targets)
#  - Latest Handback DateTime (G) is stamped with the handback time
#  - Handoff Reason (H) flips from "Ignored" to "Include" now that the row has been handed back

$wb = $excel.ActiveWorkbook

$locales = @(
    @{
        SheetName   = "zh-cn"
        MdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/c6137ed5abfebc0a68c4b9eb636fd9db6bbbee55/e2e/93a10f27-93c3-46bf-9734-25b8fdc93769.md"
        XlfUrl      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eaa8063499a34bd4fe2ca43fc6dd2da77cd89a9e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/93a10f27-93c3-46bf-9734-25b8fdc93769.48ab1842ae1b41f531cc293b10e2cc3c6c9dce8e.zh-cn.xlf"
        MdName      = "93a10f27-93c3-46bf-9734-25b8fdc93769.md"
        XlfName     = "93a10f27-93c3-46bf-9734-25b8fdc93769.48ab1842ae1b41f531cc293b10e2cc3c6c9dce8e.zh-cn.xlf"
        HandbackDT  = "2016-03-09 10:57:16"
    },
    @{
        SheetName   = "de-de"
        MdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/c6137ed5abfebc0a68c4b9eb636fd9db6bbbee55/e2e/93a10f27-93c3-46bf-9734-25b8fdc93769.md"
        XlfUrl      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3240801b428aa3c28bd001377ceafbdc07c266e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/93a10f27-93c3-46bf-9734-25b8fdc93769.48ab1842ae1b41f531cc293b10e2cc3c6c9dce8e.de-de.xlf"
        MdName      = "93a10f27-93c3-46bf-9734-25b8fdc93769.md"
        XlfName     = "93a10f27-93c3-46bf-9734-25b8fdc93769.48ab1842ae1b41f531cc293b10e2cc3c6c9dce8e.de-de.xlf"
        HandbackDT  = "2016-03-09 10:57:26"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.SheetName)

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # New "Latest Target File" (E) and "Latest Handback File" (F) hyperlinks for rows 2 & 3
    $ws.Hyperlinks.Add($ws.Range("E2"), $locale.MdUrl, [Type]::Missing, [Type]::Missing, $locale.MdName)
    $ws.Hyperlinks.Add($ws.Range("F2"), $locale.XlfUrl, [Type]::Missing, [Type]::Missing, $locale.XlfName)
    $ws.Hyperlinks.Add($ws.Range("E3"), $locale.MdUrl, [Type]::Missing, [Type]::Missing, $locale.MdName)
    $ws.Hyperlinks.Add($ws.Range("F3"), $locale.XlfUrl, [Type]::Missing, [Type]::Missing, $locale.XlfName)

    # Latest Handback DateTime (G)
    $ws.Range("G2").Value = $locale.HandbackDT
    $ws.Range("G3").Value = $locale.HandbackDT

    # Handoff Reason (H): Ignored -> Include
    $ws.Range("H2").Value = "Include"
    $ws.Range("H3").Value = "Include"
}
